$d = $word.ActiveDocument

function Replace-Text([string]$find, [string]$replace) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Replacement.ClearFormatting()
    $result = $range.Find.Execute(
        $find,
        $true,   # MatchCase
        $false,  # MatchWholeWord
        $false,  # MatchWildcards
        $false,  # MatchSoundsLike
        $false,  # MatchAllWordForms
        $true,   # Forward
        1,       # Wrap (wdFindContinue)
        $false,  # Format
        $replace,
        2        # Replace (wdReplaceAll)
    )
    if (-not $result) {
        throw "Find/Replace failed for: $find"
    }
}

Replace-Text '2023-10-06 Friday' '2023-10-07 Saturday'
Replace-Text '67×79=' '72×38='
Replace-Text '98×24=' '11×20='
Replace-Text '30×51=' '13×99='
Replace-Text '66×30=' '41×98='
Replace-Text '84×91=' '36×20='
Replace-Text '53×81=' '49×20='
Replace-Text '32×33=' '91×65='
Replace-Text '60×43=' '85×84='
Replace-Text '63×69=' '45×47='
Replace-Text '74×38=' '67×35='
Replace-Text '79×57=' '42×44='
Replace-Text '26×13=' '71×75='
Replace-Text '46×38=' '59×26='
Replace-Text '93×18=' '37×88='
Replace-Text '19×58=' '49×11='
Replace-Text '59×95=' '93×15='
Replace-Text '31×92=' '20×33='
Replace-Text '64×25=' '84×33='
Replace-Text '53×29=' '43×49='
Replace-Text '54×63=' '91×43='
Replace-Text '95×47=' '41×66='
Replace-Text '94×24=' '32×54='
Replace-Text '17×96=' '66×71='
Replace-Text '88×85=' '41×24='
Replace-Text '58×39=' '74×26='
